# Pioneer Gliders Calibration and ingest CSV
# Rename "Glider" sheet to "Moorings", update cruise/ingest identifiers from
# GL001 to GL376, add a scale-factor value, and update the UI selections
# left behind by the editing session.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename the "Glider" sheet to "Moorings" ---
$ws1.Name = "Moorings"

# Renaming the sheet turns any defined name that referenced the sheet's
# #REF! error into a bare "#REF!" (losing the sheet-name qualifier).
# Restore the explicit "Moorings!" qualifier to match the authored file.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "_FilterDatabase_0") {
        $n.RefersTo = "=Moorings!#REF!"
    }
    if ($n.Name -eq "_FilterDatabase_0_0_0") {
        $n.RefersTo = "=Moorings!#REF!"
    }
}

# --- Moorings!E2 now carries a scale-factor coefficient value ---
$ws1.Range("E2").Value = 0.0625

# --- Asset_Cal_Info: update the ingest reference designators from the
#     GL001 glider to the GL376 glider ---
$ws2.Range("A2").Value  = "CP05MOAS-GL376-01-ADCPAM000"
$ws2.Range("A3").Value  = "CP05MOAS-GL376-01-ADCPAM000"
$ws2.Range("A4").Value  = "CP05MOAS-GL376-01-ADCPAM000"
$ws2.Range("A5").Value  = "CP05MOAS-GL376-01-ADCPAM000"
$ws2.Range("A7").Value  = "CP05MOAS-GL376-02-FLORTM000"
$ws2.Range("A8").Value  = "CP05MOAS-GL376-02-FLORTM000"
$ws2.Range("A9").Value  = "CP05MOAS-GL376-02-FLORTM000"
$ws2.Range("A10").Value = "CP05MOAS-GL376-02-FLORTM000"
$ws2.Range("A12").Value = "CP05MOAS-GL376-03-CTDGVM000"
$ws2.Range("A14").Value = "CP05MOAS-GL376-04-DOSTAM000"
$ws2.Range("A16").Value = "CP05MOAS-GL376-05-PARADM000"
$ws2.Range("A18").Value = "CP05MOAS-GL376-00-ENG000000"

# --- Restore the on-screen selections left in each sheet ---
$ws2.Activate()
$ws2.Range("C28").Select()

$ws1.Activate()
$ws1.Range("D14").Select()
